# Updates cryptos list values (prices and 1h volume changes) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.454.23'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '1.803.80'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.36'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.582'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.18%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '34.82'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +5.93%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.300'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0695'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0955'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('D12').Value = '2.065.12'
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.14'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.10%  '
$ws.Range('D14').Value = '1.800.81'
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.643'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('D16').Value = '34.460.40'
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.39'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.04'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '245.28'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.46'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.16'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '173.14'
$ws.Range('D24').ClearFormats()
$ws.Range('E25').Value = '  +1.85%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.83'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +7.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.81'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.06%  '
$ws.Range('E28').Value = '  +1.87%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.01'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0531'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.84'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.46%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.25'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('E34').Value = '  -0.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.684'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.65%  '
$ws.Range('D36').Value = '1.394.57'
$ws.Range('E36').Value = '  -2.69%  '
$ws.Range('E37').Value = '  -3.86%  '
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('E39').Value = '  -1.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '83.45'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.54%  '
$ws.Range('E41').Value = '  +2.02%  '
$ws.Range('E42').Value = '  +0.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.39'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.57'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.32%  '
$ws.Range('E45').Value = '  +3.79%  '
$ws.Range('E46').Value = '  -2.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.00'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.80%  '
$ws.Range('D48').Value = '1.965.23'
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '104.83'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.31%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0131'
$ws.Range('E50').Value = '  +1.77%  '
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.17%  '
